$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-23 Friday" "2025-05-24 Saturday"

Replace-Text "18×92=" "35×62="
Replace-Text "70×66=" "74×65="
Replace-Text "70×67=" "87×38="
Replace-Text "15×13=" "18×13="
Replace-Text "34×65=" "86×74="
Replace-Text "50×75=" "69×42="
Replace-Text "62×57=" "51×71="
Replace-Text "45×71=" "49×64="
Replace-Text "55×64=" "79×65="
Replace-Text "52×76=" "31×73="
Replace-Text "70×80=" "73×24="
Replace-Text "82×59=" "71×91="
Replace-Text "90×94=" "75×20="
Replace-Text "67×37=" "48×78="
Replace-Text "48×40=" "52×22="
Replace-Text "24×44=" "12×35="
Replace-Text "34×67=" "82×96="
Replace-Text "92×14=" "45×24="
Replace-Text "75×13=" "20×19="
Replace-Text "44×19=" "14×32="
Replace-Text "13×64=" "61×26="
Replace-Text "59×46=" "99×40="
Replace-Text "73×48=" "40×34="
Replace-Text "75×39=" "83×57="
Replace-Text "90×20=" "97×17="
